$wb = $excel.ActiveWorkbook

# --- Elementary School sheet (sheet1) ---
$ws1 = $wb.Worksheets.Item(1)

# Update the text labels (boundary map change: 1/3 & 2/3 -> 1/2)
$ws1.Range("J5").Value = "1/2=>stroller"
$ws1.Range("J22").Value = "1/2->Medow Park"

# Updated student-number figures
$ws1.Range("L2").Value = 2210
$ws1.Range("M2").Value = 1105

$ws1.Range("L5").Value = 2082
$ws1.Range("M5").Value = 1041

$ws1.Range("L19").Value = 2249
$ws1.Range("M19").Value = 1125

$ws1.Range("L23").Value = 2026
$ws1.Range("M23").Value = 1013

# Add new total formula for row 26 (current boundary map total)
$ws1.Range("I26").Formula = "=SUM(I6:I8)"

# --- Selection / active-tab bookkeeping ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D3").Select()

$ws1.Activate()
$ws1.Range("C6").Select()
